$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1382.6666
$ws.Range("I111").Value = 1098.6666
$ws.Range("J111").Value = 1666.6666
$ws.Range("K111").Value = 3295.9998
$ws.Range("L111").Value = 4999.9998
$ws.Range("M111").Value = -228.9998000000001
$ws.Range("N111").Value = -11133.9998
$ws.Range("H137").Value = 1203.5
$ws.Range("I137").Value = 821.95654
$ws.Range("J137").Value = 1878.5385
$ws.Range("K137").Value = 2465.86962
$ws.Range("L137").Value = 5635.6155
$ws.Range("M137").Value = 84.13038000000006
$ws.Range("N137").Value = -10735.6155

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14978
$ws.Range("I2").Value = 606.6
$ws.Range("K2").Value = 606.6
$ws.Range("M2").Value = -493.6
$ws.Range("H32").Value = 3428.054
$ws.Range("I32").Value = 3137.1875
$ws.Range("J32").Value = 5289.6
$ws.Range("K32").Value = 3137.1875
$ws.Range("L32").Value = 5289.6
$ws.Range("M32").Value = -2850.1875
$ws.Range("N32").Value = -5863.6
$ws.Range("H43").Value = 6920.3335
$ws.Range("J43").Value = 7236
$ws.Range("L43").Value = 7236
$ws.Range("N43").Value = -7862
$ws.Range("H45").Value = 1087.069
$ws.Range("I45").Value = 1053.8235
$ws.Range("K45").Value = 1053.8235
$ws.Range("M45").Value = -676.8235
$ws.Range("H116").Value = 14978
$ws.Range("I116").Value = 606.6
$ws.Range("K116").Value = 606.6
$ws.Range("M116").Value = 1687.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14978
$ws.Range("I3").Value = 606.6
$ws.Range("K3").Value = 606.6
$ws.Range("M3").Value = -492.6
$ws.Range("H20").Value = 1266.8
$ws.Range("J20").Value = 3133.3333
$ws.Range("L20").Value = 3133.3333
$ws.Range("N20").Value = -3627.3333
$ws.Range("H94").Value = 35714816
$ws.Range("I94").Value = 35714816
$ws.Range("K94").Value = 35714816
$ws.Range("M94").Value = -35714365
$ws.Range("H105").Value = 111114550
$ws.Range("I105").Value = 111114550
$ws.Range("K105").Value = 111114550
$ws.Range("M105").Value = -111112803
$ws.Range("H107").Value = 1353.3846
$ws.Range("I107").Value = 1255.6666
$ws.Range("K107").Value = 1255.6666
$ws.Range("M107").Value = 664.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H31").Value = 1331.2354
$ws.Range("I31").Value = 1284.9574
$ws.Range("J31").Value = 1875
$ws.Range("K31").Value = 1284.9574
$ws.Range("L31").Value = 1875
$ws.Range("M31").Value = -989.9574
$ws.Range("N31").Value = -2465
$ws.Range("H34").Value = 1331.2354
$ws.Range("I34").Value = 1284.9574
$ws.Range("J34").Value = 1875
$ws.Range("K34").Value = 1284.9574
$ws.Range("L34").Value = 1875
$ws.Range("M34").Value = -1082.9574
$ws.Range("N34").Value = -2279
$ws.Range("H38").Value = 2000
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 2000
$ws.Range("N38").Value = -2754
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29407
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2422
$ws.Range("H107").Value = 426.82608
$ws.Range("I107").Value = 341.8889
$ws.Range("J107").Value = 732.6
$ws.Range("K107").Value = 341.8889
$ws.Range("L107").Value = 732.6
$ws.Range("M107").Value = 1578.1111
$ws.Range("N107").Value = -4572.6
$ws.Range("H122").Value = 1056.2
$ws.Range("I122").Value = 1023.1429
$ws.Range("J122").Value = 1133.3334
$ws.Range("K122").Value = 3069.4287
$ws.Range("L122").Value = 3400.0002
$ws.Range("M122").Value = -619.4287000000004
$ws.Range("N122").Value = -8300.0002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1289.9642
$ws.Range("I5").Value = 1332.56
$ws.Range("K5").Value = 3997.68
$ws.Range("M5").Value = -3885.68
$ws.Range("H131").Value = 12989280
$ws.Range("I131").Value = 200000400
$ws.Range("J131").Value = 2396.8472
$ws.Range("K131").Value = 600001200
$ws.Range("L131").Value = 7190.5416
$ws.Range("M131").Value = -599996160
$ws.Range("N131").Value = -17270.5416
$ws.Range("H135").Value = 1289.9642
$ws.Range("I135").Value = 1332.56
$ws.Range("K135").Value = 11993.04
$ws.Range("M135").Value = -9458.039999999999
$ws.Range("H139").Value = 1568.3572
$ws.Range("I139").Value = 1491.7693
$ws.Range("J139").Value = 1692.8125
$ws.Range("K139").Value = 4475.3079
$ws.Range("L139").Value = 5078.4375
$ws.Range("M139").Value = 664.6921000000002
$ws.Range("N139").Value = -15358.4375

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 712901.8
$ws.Range("I107").Value = 1012802.6
$ws.Range("J107").Value = 637.375
$ws.Range("K107").Value = 1012802.6
$ws.Range("L107").Value = 637.375
$ws.Range("M107").Value = -1010882.6
$ws.Range("N107").Value = -4477.375
$ws.Range("H113").Value = 1338.5333
$ws.Range("I113").Value = 1347.8572
$ws.Range("J113").Value = 1330.375
$ws.Range("K113").Value = 1347.8572
$ws.Range("L113").Value = 1330.375
$ws.Range("M113").Value = 822.1428000000001
$ws.Range("N113").Value = -5670.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1179.6
$ws.Range("I22").Value = 406.42856
$ws.Range("J22").Value = 1856.125
$ws.Range("K22").Value = 406.42856
$ws.Range("L22").Value = 1856.125
$ws.Range("M22").Value = -111.42856
$ws.Range("N22").Value = -2446.125
$ws.Range("H27").Value = 1179.6
$ws.Range("I27").Value = 406.42856
$ws.Range("J27").Value = 1856.125
$ws.Range("K27").Value = 406.42856
$ws.Range("L27").Value = 1856.125
$ws.Range("M27").Value = -299.42856
$ws.Range("N27").Value = -2070.125
$ws.Range("H40").Value = 3150.76
$ws.Range("I40").Value = 2025.2667
$ws.Range("J40").Value = 4839
$ws.Range("K40").Value = 2025.2667
$ws.Range("L40").Value = 4839
$ws.Range("M40").Value = -1889.2667
$ws.Range("N40").Value = -5111
$ws.Range("H55").Value = 233.96428
$ws.Range("J55").Value = 276.92307
$ws.Range("L55").Value = 276.92307
$ws.Range("N55").Value = -622.9230700000001
$ws.Range("H132").Value = 1764.326
$ws.Range("I132").Value = 1331.8667
$ws.Range("J132").Value = 2575.1875
$ws.Range("K132").Value = 3995.6001
$ws.Range("L132").Value = 7725.5625
$ws.Range("M132").Value = -1465.6001
$ws.Range("N132").Value = -12785.5625
$ws.Range("H136").Value = 5693.0454
$ws.Range("I136").Value = 6270.6313
$ws.Range("J136").Value = 2035
$ws.Range("K136").Value = 18811.8939
$ws.Range("L136").Value = 6105
$ws.Range("M136").Value = -16261.8939
$ws.Range("N136").Value = -11205

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 35495
$ws.Range("J75").Value = 35495
$ws.Range("L75").Value = 35495
$ws.Range("N75").Value = -37367
$ws.Range("H78").Value = 35495
$ws.Range("J78").Value = 35495
$ws.Range("L78").Value = 106485
$ws.Range("N78").Value = -115845
$ws.Range("H81").Value = 379.6
$ws.Range("I81").Value = 379.6
$ws.Range("K81").Value = 759.2
$ws.Range("M81").Value = 301.8
$ws.Range("H84").Value = 379.6
$ws.Range("I84").Value = 379.6
$ws.Range("K84").Value = 3796
$ws.Range("M84").Value = 1508
